$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "description" column (G) with a header and the same leave-request
# description repeated for every existing row.
$ws.Range("G1").Value = "description"
$ws.Range("G2").Value = "Leave request 28/03/2020"
$ws.Range("G3").Value = "Leave request 28/03/2020"
$ws.Range("G4").Value = "Leave request 28/03/2020"
$ws.Range("G5").Value = "Leave request 28/03/2020"

$ws.Range("G5").Select()
